# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Update a handful of rounded percentage/indicator figures on the Sweden
# MSME summary sheet with more precise (2-decimal) values. The source
# cells are stored as text (they already render as plain numeric-looking
# strings), so we keep writing them as text via the leading-apostrophe
# text qualifier to avoid Excel's automatic "looks like a number" type
# coercion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) - Statistical Institution source
$ws.Range("B11").Value = "'107.92"
$ws.Range("C11").Value = "'1.67"
$ws.Range("D11").Value = "'109.59"

# Employment (% of total) - Statistical Institution source
$ws.Range("B12").Value = "'33.16"
$ws.Range("C12").Value = "'30.94"
$ws.Range("D12").Value = "'64.09"

# Enterprises (% of total) - SME Associations source
$ws.Range("C33").Value = "'3.63"
$ws.Range("D33").Value = "'68.83"

# Employment (% of total) - SME Associations source
$ws.Range("B34").Value = "'25.93"
$ws.Range("C34").Value = "'39.72"
$ws.Range("D34").Value = "'65.65"

# Enterprises (% of total) - SBS Eurostat source
$ws.Range("B36").Value = "'94.58"
$ws.Range("C36").Value = "'5.26"
$ws.Range("D36").Value = "'99.85"

# Value added to the economy (% of total) - SBS Eurostat source
$ws.Range("B40").Value = "'21.69"
$ws.Range("C40").Value = "'37.81"
$ws.Range("D40").Value = "'59.51"
